# module 5: OOPs Assignment implemented with code
#
# Applies the diff:
#  - Merge the "Course Name" runs into a single run.
#  - Fix the double space after "Assignment Name -" and split the
#    assignment title into " " / "Assignment : Python - Data Structure" / " ".
#  - Split "Submission Date - August 25, 20242024" into
#    "Submission Date - " + "September 8, 2024".
#  - Replace the Github link display text with the module 3 link text,
#    collapsing all the proofErr-laden runs into a single hyperlink run.

$d = $word.ActiveDocument

function Insert-RunsXml {
    param($range, [string]$runsXml)
    $wrapper = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($wrapper)
}

function Find-ParagraphStartingWith {
    param([string]$prefix)
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Course Name - Data Science " + "with" + " Generative AI - Aug'24"
#    -> single run "Course Name - Data Science with Generative AI - Aug'24"
#    (Done via InsertXML rather than Find/Replace so the straight apostrophe
#    in "Aug'24" is not auto-corrected into a curly quote.)
# ---------------------------------------------------------------------------
$coursePara = Find-ParagraphStartingWith "Course Name"
$courseRange = $coursePara.Range
$courseStart = $courseRange.Start
$courseRange.MoveEnd(1, -1) | Out-Null
$courseRange.Delete() | Out-Null

$courseInsertPoint = $d.Range($courseStart, $courseStart)
$courseRunsXml =
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Course Name – Data Science with Generative AI - Aug' + "'" + '24</w:t></w:r>'
Insert-RunsXml $courseInsertPoint $courseRunsXml

# ---------------------------------------------------------------------------
# 2) Assignment Name paragraph:
#    "Assignment Name -  " + "Python Basics Assignment "
#    -> "Assignment Name - " + " " + "Assignment : Python - Data Structure" + " "
# ---------------------------------------------------------------------------
$assignPara = Find-ParagraphStartingWith "Assignment Name"
$assignRange = $assignPara.Range
$assignStart = $assignRange.Start
$assignRange.MoveEnd(1, -1) | Out-Null
$assignRange.Delete() | Out-Null

$assignInsertPoint = $d.Range($assignStart, $assignStart)
$assignRunsXml =
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Assignment Name – </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Assignment : Python - Data Structure</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
Insert-RunsXml $assignInsertPoint $assignRunsXml

# ---------------------------------------------------------------------------
# 3) Submission Date paragraph:
#    "Submission Date - August 25, 20242024"
#    -> "Submission Date - " + "September 8, 2024"
# ---------------------------------------------------------------------------
$subPara = $d.Paragraphs(5)
$subRange = $subPara.Range
$subRange.MoveEnd(1, -1) | Out-Null
$subRange.Delete() | Out-Null

$subPara2 = $d.Paragraphs(5)
$subInsertPoint = $d.Range($subPara2.Range.Start, $subPara2.Range.Start)
$subRunsXml =
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">Submission Date – </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>September 8, 2024</w:t></w:r>'
Insert-RunsXml $subInsertPoint $subRunsXml

# ---------------------------------------------------------------------------
# 4) Github link hyperlink display text: collapse every run + proofErr marker
#    into a single hyperlink run with the new module 3 link text.
# ---------------------------------------------------------------------------
$ghHyperlink = $d.Hyperlinks.Item(2)
$ghHyperlink.TextToDisplay = "PW_Skilles_Assignments/module_03_Python_Data_Structure_Assignment at master · Mohd-jibrail/PW_Skilles_Assignments (github.com)"
